# Updates cryptos list values (price/volume) per the source commit,
# reproducing the "Updated cryptos list ... with GitHub Actions" diff.
# All cells in columns B:E are plain text (inline strings) in the source
# workbook, so numeric-looking values are written with a leading quote
# (forces text entry) and then ClearFormats() strips the resulting
# quote-prefix / number-format styling so the cell stays styleless text,
# matching the original un-styled inlineStr cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.868.90'
$ws.Range("E2").Value = '  -1.70%  '

# Row 3
$ws.Range("D3").Value = '1.803.48'
$ws.Range("E3").Value = '  -1.19%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").Value = "'309.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.67%  '

# Row 6
$ws.Range("E6").Value = '  -0.08%  '

# Row 7
$ws.Range("D7").Value = "'0.4668"
$ws.Range("D7").ClearFormats()

# Row 8
$ws.Range("D8").Value = "'0.3702"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.98%  '

# Row 9
$ws.Range("D9").Value = "'0.07387"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.71%  '

# Row 10
$ws.Range("D10").Value = "'0.8709"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.94%  '

# Row 11
$ws.Range("E11").Value = '  -2.93%  '

# Row 12
$ws.Range("D12").Value = '1.810.41'
$ws.Range("E12").Value = '  -0.79%  '

# Row 13
$ws.Range("D13").Value = "'5.365"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.78%  '

# Row 14
$ws.Range("D14").Value = "'92.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.99%  '

# Row 15
$ws.Range("D15").Value = "'6.490"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.69%  '

# Row 16
$ws.Range("D16").Value = "'0.07028"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.45%  '

# Row 17
$ws.Range("E17").Value = '  -0.10%  '

# Row 18
$ws.Range("D18").Value = "'0.000008714"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.81%  '

# Row 19
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.04%  '

# Row 20
$ws.Range("D20").Value = "'14.69"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.89%  '

# Row 21
$ws.Range("D21").Value = '26.865.59'
$ws.Range("E21").Value = '  -1.73%  '

# Row 22
$ws.Range("D22").Value = "'5.298"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.81%  '

# Row 23
$ws.Range("E23").Value = '  -3.18%  '

# Row 24
$ws.Range("D24").Value = '2.006.48'
$ws.Range("E24").Value = '  -2.13%  '

# Row 25
$ws.Range("D25").Value = "'1.900"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.54%  '

# Row 26
$ws.Range("D26").Value = "'151.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.13%  '

# Row 27
$ws.Range("E27").Value = '  -2.02%  '

# Row 28
$ws.Range("D28").Value = "'2.133"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -8.36%  '

# Row 29
$ws.Range("D29").Value = "'5.266"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.24%  '

# Row 30
$ws.Range("D30").Value = "'115.66"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.83%  '

# Row 31
$ws.Range("D31").Value = "'0.08933"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.50%  '

# Row 32
$ws.Range("D32").Value = "'0.7597"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.47%  '

# Row 33
$ws.Range("E33").Value = '  -4.08%  '

# Row 34
$ws.Range("D34").Value = "'4.458"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.26%  '

# Row 35
$ws.Range("D35").Value = "'2.912"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.27%  '

# Row 37
$ws.Range("D37").Value = "'1.099"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.04%  '

# Row 38
$ws.Range("D38").Value = "'0.01958"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.32%  '

# Row 39
$ws.Range("D39").Value = "'0.05248"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.17%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = "'7.260"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.54%  '

# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = "'2.927"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.89%  '

# Row 42
$ws.Range("D42").Value = "'2.386"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.97%  '

# Row 43
$ws.Range("E43").Value = '  -1.07%  '

# Row 44
$ws.Range("E44").Value = '  -3.26%  '

# Row 45
$ws.Range("D45").Value = "'8.491"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.01%  '

# Row 46
$ws.Range("D46").Value = "'0.5003"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.52%  '

# Row 47
$ws.Range("D47").Value = "'10.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.28%  '

# Row 48
$ws.Range("D48").Value = "'104.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.03%  '

# Row 49
$ws.Range("E49").Value = '  -0.08%  '

# Row 50
$ws.Range("D50").Value = "'1.665"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.90%  '

# Row 51
$ws.Range("D51").Value = "'0.06293"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.85%  '
